# Apply the scraped-data refresh for the crypto price sheet (commit:
# "Updated symbol list ... with GitHub Actions"). For rows 2-51 the Price (D),
# Volume(1h) (E), and Hora (G) columns are refreshed with newly scraped values.
#
# The source cells are plain text (e.g. "319.63", "3.61%", "17") rather than
# numbers, so a leading apostrophe is used to force Excel to store the new
# value as literal text instead of auto-converting percentages/numbers; the
# style is then reset to "Normal" so the quote-prefix does not add formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Sheet, $CellRef, $Text) {
    $Sheet.Range($CellRef).Value = "'" + $Text
    $Sheet.Range($CellRef).Style = "Normal"
}

Set-TextValue $ws "D2" '319.94'
Set-TextValue $ws "E2" '3.64%'
Set-TextValue $ws "G2" '18'
Set-TextValue $ws "D3" '41.47'
Set-TextValue $ws "E3" '1.39%'
Set-TextValue $ws "G3" '18'
Set-TextValue $ws "D4" '5.242'
Set-TextValue $ws "E4" '2.23%'
Set-TextValue $ws "G4" '18'
Set-TextValue $ws "D5" '0.07729'
Set-TextValue $ws "E5" '1.32%'
Set-TextValue $ws "G5" '18'
Set-TextValue $ws "D6" '1.682'
Set-TextValue $ws "E6" '3.97%'
Set-TextValue $ws "G6" '18'
Set-TextValue $ws "D7" '0.9449'
Set-TextValue $ws "E7" '3.97%'
Set-TextValue $ws "G7" '18'
Set-TextValue $ws "E8" '-1.97%'
Set-TextValue $ws "G8" '18'
Set-TextValue $ws "D9" '0.1238'
Set-TextValue $ws "E9" '-3.68%'
Set-TextValue $ws "G9" '18'
Set-TextValue $ws "D10" '0.1835'
Set-TextValue $ws "G10" '18'
Set-TextValue $ws "D11" '0.09205'
Set-TextValue $ws "E11" '0.91%'
Set-TextValue $ws "G11" '18'
Set-TextValue $ws "D12" '0.04341'
Set-TextValue $ws "E12" '1.49%'
Set-TextValue $ws "G12" '18'
Set-TextValue $ws "D13" '0.1050'
Set-TextValue $ws "E13" '0.47%'
Set-TextValue $ws "G13" '18'
Set-TextValue $ws "D14" '0.001279'
Set-TextValue $ws "E14" '2.02%'
Set-TextValue $ws "G14" '18'
Set-TextValue $ws "D15" '0.006020'
Set-TextValue $ws "E15" '4.13%'
Set-TextValue $ws "G15" '18'
Set-TextValue $ws "E16" '-0.28%'
Set-TextValue $ws "G16" '18'
Set-TextValue $ws "D17" '4.341'
Set-TextValue $ws "E17" '1.45%'
Set-TextValue $ws "G17" '18'
Set-TextValue $ws "E18" '3.11%'
Set-TextValue $ws "G18" '18'
Set-TextValue $ws "D19" '7.659'
Set-TextValue $ws "E19" '10.94%'
Set-TextValue $ws "G19" '18'
Set-TextValue $ws "D20" '0.1354'
Set-TextValue $ws "E20" '-2.84%'
Set-TextValue $ws "G20" '18'
Set-TextValue $ws "D21" '0.2825'
Set-TextValue $ws "E21" '4.44%'
Set-TextValue $ws "G21" '18'
Set-TextValue $ws "E22" '-0.22%'
Set-TextValue $ws "G22" '18'
Set-TextValue $ws "D23" '0.001268'
Set-TextValue $ws "E23" '-0.27%'
Set-TextValue $ws "G23" '18'
Set-TextValue $ws "D24" '0.004123'
Set-TextValue $ws "E24" '1.53%'
Set-TextValue $ws "G24" '18'
Set-TextValue $ws "E25" '0.05%'
Set-TextValue $ws "G25" '18'
Set-TextValue $ws "G26" '18'
Set-TextValue $ws "G27" '18'
Set-TextValue $ws "G28" '18'
Set-TextValue $ws "G29" '18'
Set-TextValue $ws "G30" '18'
Set-TextValue $ws "G31" '18'
Set-TextValue $ws "G32" '18'
Set-TextValue $ws "G33" '18'
Set-TextValue $ws "G34" '18'
Set-TextValue $ws "G35" '18'
Set-TextValue $ws "G36" '18'
Set-TextValue $ws "G37" '18'
Set-TextValue $ws "D38" '0.02539'
Set-TextValue $ws "E38" '4.56%'
Set-TextValue $ws "G38" '18'
Set-TextValue $ws "E39" '2.17%'
Set-TextValue $ws "G39" '18'
Set-TextValue $ws "D40" '0.007777'
Set-TextValue $ws "E40" '-0.79%'
Set-TextValue $ws "G40" '18'
Set-TextValue $ws "E41" '1.23%'
Set-TextValue $ws "G41" '18'
Set-TextValue $ws "D42" '0.007367'
Set-TextValue $ws "E42" '8.29%'
Set-TextValue $ws "G42" '18'
Set-TextValue $ws "E43" '3.17%'
Set-TextValue $ws "G43" '18'
Set-TextValue $ws "D44" '0.007569'
Set-TextValue $ws "E44" '-5.98%'
Set-TextValue $ws "G44" '18'
Set-TextValue $ws "D45" '0.3459'
Set-TextValue $ws "E45" '3.72%'
Set-TextValue $ws "G45" '18'
Set-TextValue $ws "D46" '0.00006670'
Set-TextValue $ws "E46" '-3.26%'
Set-TextValue $ws "G46" '18'
Set-TextValue $ws "E47" '0.05%'
Set-TextValue $ws "G47" '18'
Set-TextValue $ws "D48" '0.1879'
Set-TextValue $ws "E48" '75.01%'
Set-TextValue $ws "G48" '18'
Set-TextValue $ws "D49" '0.004205'
Set-TextValue $ws "E49" '40.03%'
Set-TextValue $ws "G49" '18'
Set-TextValue $ws "E50" '0.05%'
Set-TextValue $ws "G50" '18'
Set-TextValue $ws "E51" '0.05%'
Set-TextValue $ws "G51" '18'
